$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "7+59=",
    "58+4=",
    "37+4=",
    "96-39=",
    "57-29=",
    "8+58=",
    "72-59=",
    "83-34=",
    "19+14=",
    "76-58=",
    "54+28=",
    "5+66=",
    "28+23=",
    "18+4=",
    "28+29=",
    "18+43=",
    "38+48=",
    "55-38=",
    "8+77=",
    "67+18=",
    "51-39=",
    "35+27=",
    "89+6=",
    "32-17=",
    "93-69=",
    "33+58=",
    "69+8=",
    "4+88=",
    "85-66=",
    "63-59=",
    "49+17=",
    "56+15=",
    "95-56=",
    "12+29=",
    "94-6=",
    "79+2=",
    "18+46=",
    "83-49=",
    "81-55=",
    "28+4=",
    "61-49=",
    "43+29=",
    "65-28=",
    "15+18=",
    "57-19=",
    "23+9=",
    "51-8=",
    "96-87=",
    "67+18=",
    "29+36=",
    "26+47=",
    "83-26=",
    "64-48=",
    "29+65=",
    "75-68=",
    "70-22=",
    "28+4=",
    "50-32=",
    "38+47=",
    "37+9=",
    "4+7=",
    "57-8=",
    "22+49=",
    "94-17=",
    "39+3=",
    "83-27=",
    "65-38=",
    "71-8=",
    "53-34=",
    "37+9=",
    "88-39=",
    "29+43=",
    "24+59=",
    "53-15=",
    "7+46=",
    "89+6=",
    "26+15=",
    "25+17=",
    "39+7=",
    "90-24=",
    "39+23=",
    "21-18=",
    "51-7=",
    "62-54=",
    "73-26=",
    "13+8=",
    "63-28=",
    "36+29=",
    "80-42=",
    "19+34=",
    "40-38=",
    "15+47=",
    "72-59=",
    "7+29=",
    "9+86=",
    "66-47=",
    "63-46=",
    "52-28=",
    "13-7=",
    "68-29="
)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        if ($idx -lt $values.Length) {
            $cell = $t.Cell($r, $c)
            $cell.Range.Text = $values[$idx]
        }
        $idx = $idx + 1
    }
}
Write-Host "Updated" $idx "cells"
